# Add support for Benchling URLs as inputs: insert a new "BenchlingUrlSource"
# worksheet right after "AddGeneIdSource" (and before "GenomeCoordinatesSource"),
# with the same header/validation layout as the other *Source sheets.

$wb = $excel.ActiveWorkbook
$originallyActive = $wb.ActiveSheet

$anchor = $wb.Worksheets.Item("GenomeCoordinatesSource")
$new = $wb.Worksheets.Add($anchor)
$new.Name = "BenchlingUrlSource"

$headers = @("repository_name", "repository_id", "input", "output", "type", "output_name", "id")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $new.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$validationRange = $new.Range("A2:A1048576")
$validationRange.Validation.Add(3, 1, 1, '"addgene,genbank"')
$validationRange.Validation.ShowDropDown = $false
$validationRange.Validation.ShowInput = $false
$validationRange.Validation.ShowError = $false

# Restore whatever sheet was active before this edit (adding a sheet makes
# it the active one, which isn't part of the intended change).
$originallyActive.Activate()
